$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply text number format to columns A and B (header row 1, data rows 2-40, footer rows 41-42)
# so that the new dotted "sheet.item" multileader labels are stored literally as text
$ws.Range("A1:B42").NumberFormat = "@"

# Replace the plain numeric position codes in column A with text-based multileader labels
$ws.Range("A2").Value = "6.2"
$ws.Range("A3").Value = "6.2"
$ws.Range("A4").Value = "8.8"
$ws.Range("A5").Value = "1.16"
$ws.Range("A6").Value = "6.7"
$ws.Range("A7").Value = "2.42"
$ws.Range("A8").Value = "8.24"
$ws.Range("A9").Value = "1.11"
$ws.Range("A10").Value = "2.19"
$ws.Range("A11").Value = "4.49"
$ws.Range("A12").Value = "8.6"
$ws.Range("A13").Value = "2.32"
$ws.Range("A14").Value = "8.49"
$ws.Range("A15").Value = "9.46"
$ws.Range("A16").Value = "8.10"
$ws.Range("A17").Value = "6.25"
$ws.Range("A18").Value = "10.2"
$ws.Range("A19").Value = "7.6"
$ws.Range("A20").Value = "5.4"
$ws.Range("A21").Value = "10.30"
$ws.Range("A22").Value = "7.32"
$ws.Range("A23").Value = "4.13"
$ws.Range("A24").Value = "3.3"
$ws.Range("A25").Value = "6.16"
$ws.Range("A26").Value = "10.47"
$ws.Range("A27").Value = "9.20"
$ws.Range("A28").Value = "1.32"
$ws.Range("A29").Value = "1.29"
$ws.Range("A30").Value = "5.41"
$ws.Range("A31").Value = "8.41"
$ws.Range("A32").Value = "3.22"
$ws.Range("A33").Value = "4.32"
$ws.Range("A34").Value = "1.48"
$ws.Range("A35").Value = "10.27"
$ws.Range("A36").Value = "7.4"
$ws.Range("A37").Value = "4.37"
$ws.Range("A38").Value = "9.34"
$ws.Range("A39").Value = "2.48"
$ws.Range("A40").Value = "10.8"
$ws.Range("A41").Value = "8.6"
